$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 98, shifting existing data (rows 98-170) down to rows 100-172
$ws.Rows("98:99").Insert()

# Populate new row 98 with the latest week of data (2022-01-05 / serial 44566)
$ws.Cells.Item(98, 1).Value2 = 11
$ws.Cells.Item(98, 2).Value2 = 'Vega Monumental Concepción'
$ws.Cells.Item(98, 3).Value2 = 'Bíobío'
$ws.Cells.Item(98, 4).Value2 = 44566
$ws.Cells.Item(98, 5).Value2 = 8
$ws.Cells.Item(98, 6).Value2 = 100114013
$ws.Cells.Item(98, 7).Value2 = 'Zanahoria'
$ws.Cells.Item(98, 8).Value2 = 'Sin especificar'
$ws.Cells.Item(98, 9).Value2 = 'Primera'
$ws.Cells.Item(98, 10).Value2 = 250
$ws.Cells.Item(98, 11).Value2 = 7500
$ws.Cells.Item(98, 12).Value2 = 8000
$ws.Cells.Item(98, 13).Value2 = 7700
$ws.Cells.Item(98, 14).Value2 = '$/saco 20 kilos'
$ws.Cells.Item(98, 15).Value2 = 'Chillán'
$ws.Cells.Item(98, 16).Value2 = 385
$ws.Cells.Item(98, 17).Value2 = 20
$ws.Cells.Item(98, 18).Value2 = 'Hortaliza'

# Populate new row 99 with the latest week of data (2022-01-05 / serial 44566)
$ws.Cells.Item(99, 1).Value2 = 11
$ws.Cells.Item(99, 2).Value2 = 'Vega Monumental Concepción'
$ws.Cells.Item(99, 3).Value2 = 'Bíobío'
$ws.Cells.Item(99, 4).Value2 = 44566
$ws.Cells.Item(99, 5).Value2 = 8
$ws.Cells.Item(99, 6).Value2 = 100114013
$ws.Cells.Item(99, 7).Value2 = 'Zanahoria'
$ws.Cells.Item(99, 8).Value2 = 'Sin especificar'
$ws.Cells.Item(99, 9).Value2 = 'Primera'
$ws.Cells.Item(99, 10).Value2 = 250
$ws.Cells.Item(99, 11).Value2 = 7500
$ws.Cells.Item(99, 12).Value2 = 8000
$ws.Cells.Item(99, 13).Value2 = 7700
$ws.Cells.Item(99, 14).Value2 = '$/saco 20 kilos'
$ws.Cells.Item(99, 15).Value2 = 'Región Metropolitana'
$ws.Cells.Item(99, 16).Value2 = 385
$ws.Cells.Item(99, 17).Value2 = 20
$ws.Cells.Item(99, 18).Value2 = 'Hortaliza'
